# Update "想去人数" (column F) counts across the four sheets of the workbook.
# Values were incremented on a later data refresh (gh-pages generated output).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2825
$ws1.Range("F3").Value = 1142
$ws1.Range("F4").Value = 20639
$ws1.Range("F6").Value = 2653
$ws1.Range("F9").Value = 494
$ws1.Range("F12").Value = 259
$ws1.Range("F17").Value = 246
$ws1.Range("F18").Value = 5
$ws1.Range("F19").Value = 406
$ws1.Range("F20").Value = 17

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 130

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6096
$ws3.Range("F3").Value = 685
$ws3.Range("F5").Value = 1472
$ws3.Range("F6").Value = 46

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6096
$ws4.Range("F3").Value = 685
$ws4.Range("F5").Value = 1472
$ws4.Range("F6").Value = 2825
$ws4.Range("F7").Value = 1142
$ws4.Range("F8").Value = 20639
$ws4.Range("F14").Value = 2653
$ws4.Range("F17").Value = 46
$ws4.Range("F19").Value = 494
$ws4.Range("F22").Value = 259
$ws4.Range("F34").Value = 246
$ws4.Range("F35").Value = 130
$ws4.Range("F36").Value = 130
$ws4.Range("F37").Value = 5
$ws4.Range("F38").Value = 406
$ws4.Range("F40").Value = 17
